# Apply ranking reshuffle to "max-arrecad" and "tx-sucesso" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: max-arrecad ---
$wsMax = $wb.Worksheets.Item("max-arrecad")

$wsMax.Range("A2").Value = "angelo_agostini"
$wsMax.Range("A3").Value = "humor"
$wsMax.Range("A4").Value = "religiosidade"
$wsMax.Range("A5").Value = "terror"
$wsMax.Range("A6").Value = "hqmix"
$wsMax.Range("A7").Value = "jogos"

$wsMax.Range("A17").Value = "zine"
$wsMax.Range("A18").Value = "herois"

# --- Sheet: tx-sucesso ---
$wsTx = $wb.Worksheets.Item("tx-sucesso")

$wsTx.Range("A15").Value = "erotismo"
$wsTx.Range("A16").Value = "politica"
